$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: new data row (CityA). A3 needs the "darker" font style currently
# used by the (soon to be removed) blank A6/A7 cells, so grab that format
# via copy/paste-special before it disappears.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "CityA"
$ws.Range("A6").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B3").Value = 157113
$ws.Range("C3").Value = 1.8
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "2021/1/1,2021/1/11,2021/2/11,2021/2/23,2021/3/20,2021/4/29,2021/5/3,2021/5/4,2021/5/5,2021/7/20,2021/8/8,2021/8/9,2021/9/20,2021/9/23,2021/11/3,2021/11/23,2022/1/10,2022/2/11,2022/2/23,2022/3/21,2022/4/29,2022/5/3,2022/5/4,2022/5/5,2022/7/18,2022/8/11,2022/9/19,2022/9/23,2022/10/10,2022/11/3,2022/11/23"

# ---------------------------------------------------------------------------
# Row 4: new data row (Kyoto_kyotoshi). C4/D4 must hold plain numbers (not
# text), so first switch those two cells away from the inherited "text"
# number format (copied from the General-formatted B2 cell) before writing
# the numeric values.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Kyoto_kyotoshi"
$ws.Range("B4").Value = 778000

$ws.Range("B2").Copy()
$ws.Range("C4:D4").PasteSpecial(-4122)   # xlPasteFormats -> General format
$ws.Range("C4").Value = 1.8
$ws.Range("D4").Value = 0

$ws.Range("E4").Value = "2021/1/1,2021/1/11,2021/2/11,2021/2/23,2021/3/20,2021/4/29,2021/5/3,2021/5/4,2021/5/5,2021/7/20,2021/8/8,2021/8/9,2021/9/20,2021/9/23,2021/11/3,2021/11/23,2022/1/10,2022/2/11,2022/2/23,2022/3/21,2022/4/29,2022/5/3,2022/5/4,2022/5/5,2022/7/18,2022/8/11,2022/9/19,2022/9/23,2022/10/10,2022/11/3,2022/11/23"

# ---------------------------------------------------------------------------
# Row 5: new data row (Kanagawa). D5 keeps the inherited "text" number
# format but must still store a genuine number (0) rather than a text
# value, matching the source workbook. To get a numeric value into a
# text-formatted cell we temporarily switch it to General, assign the
# number, then restore the original text format (its value is left alone
# by a format-only paste).
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Kanagawa"
$ws.Range("B5").Value = 1241200
$ws.Range("C5").Value = 1.8

$ws.Range("B2").Copy()
$ws.Range("D5").PasteSpecial(-4122)   # xlPasteFormats -> General format
$ws.Range("D5").Value = 0
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)   # xlPasteFormats -> restore text format

$ws.Range("E5").Value = "2021/1/1,2021/1/11,2021/2/11,2021/2/23,2021/3/20,2021/4/29,2021/5/3,2021/5/4,2021/5/5,2021/7/20,2021/8/8,2021/8/9,2021/9/20,2021/9/23,2021/11/3,2021/11/23,2022/1/10,2022/2/11,2022/2/23,2022/3/21,2022/4/29,2022/5/3,2022/5/4,2022/5/5,2022/7/18,2022/8/11,2022/9/19,2022/9/23,2022/10/10,2022/11/3,2022/11/23"

# ---------------------------------------------------------------------------
# The former rows 6-10 are no longer part of the table; drop them so the
# sheet dimension shrinks back down to A1:E5.
# ---------------------------------------------------------------------------
$ws.Rows("6:10").Delete()

# ---------------------------------------------------------------------------
# Selection moves to A4 (single cell) rather than the old A3:O10 block.
# ---------------------------------------------------------------------------
$ws.Range("A4").Select() | Out-Null
